$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J, matching style of existing header row (style index 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I and J columns, rows 2-51
$values = @(
    @(2, 10, 10),
    @(3, 7, 7),
    @(4, 9, 9),
    @(5, 8, 8),
    @(6, 6, 6),
    @(7, 6, 7),
    @(8, 6, 7),
    @(9, 8, 8),
    @(10, 6, 6),
    @(11, 8, 8),
    @(12, 5, 5),
    @(13, 6, 7),
    @(14, 10, 10),
    @(15, 5, 6),
    @(16, 9, 9),
    @(17, 9, 9),
    @(18, 9, 9),
    @(19, 5, 5),
    @(20, 5, 5),
    @(21, 8, 9),
    @(22, 3, 3),
    @(23, 8, 8),
    @(24, 6, 6),
    @(25, 7, 7),
    @(26, 6, 6),
    @(27, 8, 8),
    @(28, 9, 9),
    @(29, 5, 5),
    @(30, 7, 7),
    @(31, 6, 7),
    @(32, 7, 7),
    @(33, 5, 5),
    @(34, 8, 8),
    @(35, 8, 8),
    @(36, 10, 10),
    @(37, 8, 9),
    @(38, 9, 9),
    @(39, 6, 7),
    @(40, 10, 10),
    @(41, 8, 8),
    @(42, 8, 8),
    @(43, 8, 8),
    @(44, 7, 7),
    @(45, 8, 8),
    @(46, 8, 8),
    @(47, 7, 8),
    @(48, 7, 7),
    @(49, 6, 6),
    @(50, 8, 8),
    @(51, 6, 6)
)

foreach ($row in $values) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
